# Refresh the cryptos price/volume table with the latest scrape:
# updates Price (D) and Volume(1h) (E) values for most rows, and
# corrects the coin order for two re-ranked pairs (rows 29-30 and
# 42-44) where Coin/Link/Price/Volume all moved to a different row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.501.07"
$ws.Range("E2").Value = "  -2.82%  "

$ws.Range("D3").Value = "2.457.97"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.17"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.53"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -6.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.64%  "

$ws.Range("D9").Value = "2.466.99"
$ws.Range("E9").Value = "  -4.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -5.18%  "

$ws.Range("D14").Value = "2.894.77"
$ws.Range("E14").Value = "  -3.62%  "

$ws.Range("D15").Value = "58.416.69"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.71"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -5.48%  "

$ws.Range("E17").Value = "  -3.87%  "

$ws.Range("D18").Value = "2.463.57"
$ws.Range("E18").Value = "  -4.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.82"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("E20").Value = "  -3.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.06"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -3.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.58"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -5.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.46"
$ws.Range("D28").NumberFormat = "General"

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.54"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -7.39%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0754"
$ws.Range("E30").Value = "  -5.32%  "

$ws.Range("E31").Value = "  -3.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.68"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.12"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -5.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.26"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.54%  "

$ws.Range("E36").Value = "  -8.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.04"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -8.08%  "

$ws.Range("E38").Value = "  -6.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.47"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.801"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("E41").Value = "  -4.64%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.09"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -9.54%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "274.74"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -8.98%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.84"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.588"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.82%  "

$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.89"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -4.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0505"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0218"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -5.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.16"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -6.11%  "
